$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 13.83295502832081
$ws.Cells.Item(2, 4).Value = 9.029302741564841
$ws.Cells.Item(2, 5).Value = 14.52394921611808
$ws.Cells.Item(2, 6).Value = 35.45959808072185
$ws.Cells.Item(2, 7).Value = 3.673527352855625
$ws.Cells.Item(2, 9).Value = 28.29257046083297
$ws.Cells.Item(2, 10).Value = 10.92615328302206
$ws.Cells.Item(2, 11).Value = 9.828957127267065
$ws.Cells.Item(2, 12).Value = 10.22711312151314
$ws.Cells.Item(2, 13).Value = 15.1612007403791
$ws.Cells.Item(2, 14).Value = 21.19685984026296
$ws.Cells.Item(2, 15).Value = 26.85140851968896
$ws.Cells.Item(3, 2).Value = 13.73190790567448
$ws.Cells.Item(3, 4).Value = 9.036440917061068
$ws.Cells.Item(3, 5).Value = 14.54924177700242
$ws.Cells.Item(3, 6).Value = 35.53472361598589
$ws.Cells.Item(3, 7).Value = 3.675274722463056
$ws.Cells.Item(3, 9).Value = 28.39013254910793
$ws.Cells.Item(3, 10).Value = 10.93712547443855
$ws.Cells.Item(3, 11).Value = 9.616825731504202
$ws.Cells.Item(3, 12).Value = 10.21213756028255
$ws.Cells.Item(3, 13).Value = 15.13785495415025
$ws.Cells.Item(3, 14).Value = 21.2542231130258
$ws.Cells.Item(3, 15).Value = 26.91901384590175
$ws.Cells.Item(4, 2).Value = 13.67165181720267
$ws.Cells.Item(4, 4).Value = 9.041710670109355
$ws.Cells.Item(4, 5).Value = 14.56572847883234
$ws.Cells.Item(4, 6).Value = 35.58749559921725
$ws.Cells.Item(4, 7).Value = 3.676405822599762
$ws.Cells.Item(4, 9).Value = 28.45382280442382
$ws.Cells.Item(4, 10).Value = 10.9442497856249
$ws.Cells.Item(4, 11).Value = 9.485617567018375
$ws.Cells.Item(4, 12).Value = 10.20419786411763
$ws.Cells.Item(4, 13).Value = 15.12532734825078
$ws.Cells.Item(4, 14).Value = 21.29113330254654
$ws.Cells.Item(4, 15).Value = 26.96525736011742
$ws.Cells.Item(5, 2).Value = 13.64756803470552
$ws.Cells.Item(5, 4).Value = 9.044081727047052
$ws.Cells.Item(5, 5).Value = 14.57268817304919
$ws.Cells.Item(5, 6).Value = 35.61067016915177
$ws.Cells.Item(5, 7).Value = 3.67688143652292
$ws.Cells.Item(5, 9).Value = 28.48073050152425
$ws.Cells.Item(5, 10).Value = 10.94725065884176
$ws.Cells.Item(5, 11).Value = 9.431982741765953
$ws.Cells.Item(5, 12).Value = 10.20128053669101
$ws.Cells.Item(5, 13).Value = 15.12068052877473
$ws.Cells.Item(5, 14).Value = 21.30660039928879
$ws.Cells.Item(5, 15).Value = 26.98529135887691
$ws.Cells.Item(6, 2).Value = 13.64359800977278
$ws.Cells.Item(6, 4).Value = 9.044488959913144
$ws.Cells.Item(6, 5).Value = 14.57385841208157
$ws.Cells.Item(6, 6).Value = 35.61461908378153
$ws.Cells.Item(6, 7).Value = 3.676961299885624
$ws.Cells.Item(6, 9).Value = 28.48525612270984
$ws.Cells.Item(6, 10).Value = 10.94775485750336
$ws.Cells.Item(6, 11).Value = 9.423068895880188
$ws.Cells.Item(6, 12).Value = 10.20081540250624
$ws.Cells.Item(6, 13).Value = 15.11993671757403
$ws.Cells.Item(6, 14).Value = 21.30919445740331
$ws.Cells.Item(6, 15).Value = 26.9886897971171
$ws.Cells.Item(7, 2).Value = 13.67132507953527
$ws.Cells.Item(7, 4).Value = 9.041741740951675
$ws.Cells.Item(7, 5).Value = 14.56582136217594
$ws.Cells.Item(7, 6).Value = 35.58780138147187
$ws.Cells.Item(7, 7).Value = 3.676412177391334
$ws.Cells.Item(7, 9).Value = 28.45418182915046
$ws.Cells.Item(7, 10).Value = 10.94428986068787
$ws.Cells.Item(7, 11).Value = 9.484894802506384
$ws.Cells.Item(7, 12).Value = 10.20415722862133
$ws.Cells.Item(7, 13).Value = 15.12526281901969
$ws.Cells.Item(7, 14).Value = 21.2913401711901
$ws.Cells.Item(7, 15).Value = 26.96552273110282
$ws.Cells.Item(8, 2).Value = 13.79775623940875
$ws.Cells.Item(8, 4).Value = 9.031580262648143
$ws.Cells.Item(8, 5).Value = 14.53247182903916
$ws.Cells.Item(8, 6).Value = 35.48412157861734
$ws.Cells.Item(8, 7).Value = 3.674117791634776
$ws.Cells.Item(8, 9).Value = 28.325424649498
$ws.Cells.Item(8, 10).Value = 10.92985625765289
$ws.Cells.Item(8, 11).Value = 9.756057165469006
$ws.Cells.Item(8, 12).Value = 10.22169054113851
$ws.Cells.Item(8, 13).Value = 15.15277854766332
$ws.Cells.Item(8, 14).Value = 21.21628889225676
$ws.Cells.Item(8, 15).Value = 26.87373586852126
$ws.Cells.Item(9, 2).Value = 14.05888569790821
$ws.Cells.Item(9, 4).Value = 9.018665034729356
$ws.Cells.Item(9, 5).Value = 14.4746398256394
$ws.Cells.Item(9, 6).Value = 35.33357722746194
$ws.Cells.Item(9, 7).Value = 3.670078333905328
$ws.Cells.Item(9, 9).Value = 28.10292396586171
$ws.Cells.Item(9, 10).Value = 10.90461363215726
$ws.Cells.Item(9, 11).Value = 10.27705905597197
$ws.Cells.Item(9, 12).Value = 10.26591857927906
$ws.Cells.Item(9, 13).Value = 15.22089105012549
$ws.Cells.Item(9, 14).Value = 21.08245844834338
$ws.Cells.Item(9, 15).Value = 26.73133856963921
$ws.Cells.Item(10, 2).Value = 14.25742883518287
$ws.Cells.Item(10, 4).Value = 9.013415581612126
$ws.Cells.Item(10, 5).Value = 14.43672616440555
$ws.Cells.Item(10, 6).Value = 35.25519902378308
$ws.Cells.Item(10, 7).Value = 3.6673880283178
$ws.Cells.Item(10, 9).Value = 27.95766078628438
$ws.Cells.Item(10, 10).Value = 10.88791793087864
$ws.Cells.Item(10, 11).Value = 10.6493079496198
$ws.Cells.Item(10, 12).Value = 10.3042515151134
$ws.Cells.Item(10, 13).Value = 15.27930591062698
$ws.Cells.Item(10, 14).Value = 20.99218924742073
$ws.Cells.Item(10, 15).Value = 26.64968645293398
$ws.Cells.Item(11, 2).Value = 14.34890721408077
$ws.Cells.Item(11, 4).Value = 9.011940093155081
$ws.Cells.Item(11, 5).Value = 14.42046393459394
$ws.Cells.Item(11, 6).Value = 35.22654688977462
$ws.Cells.Item(11, 7).Value = 3.666223788740454
$ws.Cells.Item(11, 9).Value = 27.89551486030535
$ws.Cells.Item(11, 10).Value = 10.8807208710867
$ws.Cells.Item(11, 11).Value = 10.81559702196956
$ws.Cells.Item(11, 12).Value = 10.32291980235561
$ws.Cells.Item(11, 13).Value = 15.30763898545591
$ws.Cells.Item(11, 14).Value = 20.95285582556063
$ws.Cells.Item(11, 15).Value = 26.61753590067325
$ws.Cells.Item(12, 2).Value = 14.38368916227374
$ws.Cells.Item(12, 4).Value = 9.01151186240171
$ws.Cells.Item(12, 5).Value = 14.41444688048059
$ws.Cells.Item(12, 6).Value = 35.21670399546806
$ws.Cells.Item(12, 7).Value = 3.665791445408932
$ws.Cells.Item(12, 9).Value = 27.87254658086765
$ws.Cells.Item(12, 10).Value = 10.87805248678856
$ws.Cells.Item(12, 11).Value = 10.87806911525004
$ws.Cells.Item(12, 12).Value = 10.33016241055252
$ws.Cells.Item(12, 13).Value = 15.31861560559364
$ws.Cells.Item(12, 14).Value = 20.93820887813637
$ws.Cells.Item(12, 15).Value = 26.6060798229
$ws.Cells.Item(13, 2).Value = 14.37619238404795
$ws.Cells.Item(13, 4).Value = 9.011598296992549
$ws.Cells.Item(13, 5).Value = 14.41573649363243
$ws.Cells.Item(13, 6).Value = 35.21877905213152
$ws.Cells.Item(13, 7).Value = 3.665884179597535
$ws.Cells.Item(13, 9).Value = 27.87746809346844
$ws.Cells.Item(13, 10).Value = 10.87862464028074
$ws.Cells.Item(13, 11).Value = 10.86463768394318
$ws.Cells.Item(13, 12).Value = 10.32859493530363
$ws.Cells.Item(13, 13).Value = 15.31624067744583
$ws.Cells.Item(13, 14).Value = 20.94135235790721
$ws.Cells.Item(13, 15).Value = 26.60851512458403
$ws.Cells.Item(14, 2).Value = 14.35176606639428
$ws.Cells.Item(14, 4).Value = 9.011902251032728
$ws.Cells.Item(14, 5).Value = 14.41996608308267
$ws.Cells.Item(14, 6).Value = 35.22571692682438
$ws.Cells.Item(14, 7).Value = 3.666188048893507
$ws.Cells.Item(14, 9).Value = 27.89361392658813
$ws.Cells.Item(14, 10).Value = 10.88050020073583
$ws.Cells.Item(14, 11).Value = 10.82074694520467
$ws.Cells.Item(14, 12).Value = 10.32351220343749
$ws.Cells.Item(14, 13).Value = 15.30853711101798
$ws.Cells.Item(14, 14).Value = 20.95164585291288
$ws.Cells.Item(14, 15).Value = 26.61657899675818
$ws.Cells.Item(15, 2).Value = 14.33682184748818
$ws.Cells.Item(15, 4).Value = 9.012105405561249
$ws.Cells.Item(15, 5).Value = 14.42257518853346
$ws.Cells.Item(15, 6).Value = 35.23009771993151
$ws.Cells.Item(15, 7).Value = 3.66637528713042
$ws.Cells.Item(15, 9).Value = 27.9035772759187
$ws.Cells.Item(15, 10).Value = 10.88165645011084
$ws.Cells.Item(15, 11).Value = 10.79379605290461
$ws.Cells.Item(15, 12).Value = 10.32042134685197
$ws.Cells.Item(15, 13).Value = 15.30385051687476
$ws.Cells.Item(15, 14).Value = 20.95798315047889
$ws.Cells.Item(15, 15).Value = 26.62161195454921
$ws.Cells.Item(16, 2).Value = 14.25147165343859
$ws.Cells.Item(16, 4).Value = 9.013530310003816
$ws.Cells.Item(16, 5).Value = 14.43780871270191
$ws.Cells.Item(16, 6).Value = 35.2572124319673
$ws.Cells.Item(16, 7).Value = 3.6674653098153
$ws.Cells.Item(16, 9).Value = 27.96180128086777
$ws.Cells.Item(16, 10).Value = 10.88839626257613
$ws.Cells.Item(16, 11).Value = 10.63837403303905
$ws.Cells.Item(16, 12).Value = 10.30305595390151
$ws.Cells.Item(16, 13).Value = 15.2774892197951
$ws.Cells.Item(16, 14).Value = 20.99479450864563
$ws.Cells.Item(16, 15).Value = 26.65188810098888
$ws.Cells.Item(17, 2).Value = 14.1993909209368
$ws.Cells.Item(17, 4).Value = 9.014637660601942
$ws.Cells.Item(17, 5).Value = 14.44740586320246
$ws.Cells.Item(17, 6).Value = 35.27564010931281
$ws.Cells.Item(17, 7).Value = 3.668149237856742
$ws.Cells.Item(17, 9).Value = 27.99852708480036
$ws.Cells.Item(17, 10).Value = 10.89263266529517
$ws.Cells.Item(17, 11).Value = 10.54220310779721
$ws.Cells.Item(17, 12).Value = 10.29271548059453
$ws.Cells.Item(17, 13).Value = 15.26176427223245
$ws.Cells.Item(17, 14).Value = 21.01781953508386
$ws.Cells.Item(17, 15).Value = 26.67174093165744
$ws.Cells.Item(18, 2).Value = 14.16954620732965
$ws.Cells.Item(18, 4).Value = 9.015360514018774
$ws.Cells.Item(18, 5).Value = 14.45301862350285
$ws.Cells.Item(18, 6).Value = 35.28689834530853
$ws.Cells.Item(18, 7).Value = 3.668548226845237
$ws.Cells.Item(18, 9).Value = 28.02002122586234
$ws.Cells.Item(18, 10).Value = 10.89510679787044
$ws.Cells.Item(18, 11).Value = 10.4866044684046
$ws.Cells.Item(18, 12).Value = 10.28688395474272
$ws.Cells.Item(18, 13).Value = 15.25288569246689
$ws.Cells.Item(18, 14).Value = 21.0312258668915
$ws.Cells.Item(18, 15).Value = 26.68362973060491
$ws.Cells.Item(19, 2).Value = 14.15946110370637
$ws.Cells.Item(19, 4).Value = 9.015620040399041
$ws.Cells.Item(19, 5).Value = 14.4549349508096
$ws.Cells.Item(19, 6).Value = 35.29082338745622
$ws.Cells.Item(19, 7).Value = 3.668684282714274
$ws.Cells.Item(19, 9).Value = 28.02736241981515
$ws.Cells.Item(19, 10).Value = 10.89595093851374
$ws.Cells.Item(19, 11).Value = 10.46773288627421
$ws.Cells.Item(19, 12).Value = 10.28492953760762
$ws.Cells.Item(19, 13).Value = 15.24990823347096
$ws.Cells.Item(19, 14).Value = 21.03579303954563
$ws.Cells.Item(19, 15).Value = 26.68773577219453
$ws.Cells.Item(20, 2).Value = 14.204923728222
$ws.Cells.Item(20, 4).Value = 9.014510892416443
$ws.Cells.Item(20, 5).Value = 14.44637463644687
$ws.Cells.Item(20, 6).Value = 35.27361023307883
$ws.Cells.Item(20, 7).Value = 3.668075852058497
$ws.Cells.Item(20, 9).Value = 27.9945792281736
$ws.Cells.Item(20, 10).Value = 10.89217781711099
$ws.Cells.Item(20, 11).Value = 10.55247046985225
$ws.Cells.Item(20, 12).Value = 10.29380425743
$ws.Cells.Item(20, 13).Value = 15.26342107957913
$ws.Cells.Item(20, 14).Value = 21.01535162528543
$ws.Cells.Item(20, 15).Value = 26.66957891847353
$ws.Cells.Item(21, 2).Value = 14.35893704497433
$ws.Cells.Item(21, 4).Value = 9.011809436331976
$ws.Cells.Item(21, 5).Value = 14.41871992490399
$ws.Cells.Item(21, 6).Value = 35.22365177473768
$ws.Cells.Item(21, 7).Value = 3.666098563964639
$ws.Cells.Item(21, 9).Value = 27.88885617552967
$ws.Cells.Item(21, 10).Value = 10.87994775841296
$ws.Cells.Item(21, 11).Value = 10.83365268803435
$ws.Cells.Item(21, 12).Value = 10.32500045028205
$ws.Cells.Item(21, 13).Value = 15.31079316361629
$ws.Cells.Item(21, 14).Value = 20.94861568814428
$ws.Cells.Item(21, 15).Value = 26.61419093433157
$ws.Cells.Item(22, 2).Value = 14.46040394523988
$ws.Cells.Item(22, 4).Value = 9.010804208807686
$ws.Cells.Item(22, 5).Value = 14.4014681670043
$ws.Cells.Item(22, 6).Value = 35.1968706876541
$ws.Cells.Item(22, 7).Value = 3.664855986422603
$ws.Cells.Item(22, 9).Value = 27.82305314099645
$ws.Cells.Item(22, 10).Value = 10.87228676062731
$ws.Cells.Item(22, 11).Value = 11.01449426066945
$ws.Cells.Item(22, 12).Value = 10.34639750182378
$ws.Cells.Item(22, 13).Value = 15.34319343834105
$ws.Cells.Item(22, 14).Value = 20.90644357227813
$ws.Cells.Item(22, 15).Value = 26.58218063881876
$ws.Cells.Item(23, 2).Value = 14.40618353164785
$ws.Cells.Item(23, 4).Value = 9.011271390474519
$ws.Cells.Item(23, 5).Value = 14.41060069640669
$ws.Cells.Item(23, 6).Value = 35.21062723712166
$ws.Cells.Item(23, 7).Value = 3.665514639713864
$ws.Cells.Item(23, 9).Value = 27.85787240328179
$ws.Cells.Item(23, 10).Value = 10.87634527231929
$ws.Cells.Item(23, 11).Value = 10.918262164042
$ws.Cells.Item(23, 12).Value = 10.33488644229742
$ws.Cells.Item(23, 13).Value = 15.32577094911735
$ws.Cells.Item(23, 14).Value = 20.92881988049293
$ws.Cells.Item(23, 15).Value = 26.59888168581204
$ws.Cells.Item(24, 2).Value = 14.2024220413439
$ws.Cells.Item(24, 4).Value = 9.014567935699734
$ws.Cells.Item(24, 5).Value = 14.44684055740348
$ws.Cells.Item(24, 6).Value = 35.27452587207969
$ws.Cells.Item(24, 7).Value = 3.668109011740916
$ws.Cells.Item(24, 9).Value = 27.99636287045945
$ws.Cells.Item(24, 10).Value = 10.89238333382542
$ws.Cells.Item(24, 11).Value = 10.54782955239953
$ws.Cells.Item(24, 12).Value = 10.29331166807654
$ws.Cells.Item(24, 13).Value = 15.26267153225045
$ws.Cells.Item(24, 14).Value = 21.01646684098922
$ws.Cells.Item(24, 15).Value = 26.67055488463982
$ws.Cells.Item(25, 2).Value = 13.98696826320061
$ws.Cells.Item(25, 4).Value = 9.021411834761617
$ws.Cells.Item(25, 5).Value = 14.48947874396593
$ws.Cells.Item(25, 6).Value = 35.36864701960896
$ws.Cells.Item(25, 7).Value = 3.671122181415791
$ws.Cells.Item(25, 9).Value = 28.15991340799543
$ws.Cells.Item(25, 10).Value = 10.9111163659411
$ws.Cells.Item(25, 11).Value = 10.13768938658819
$ws.Cells.Item(25, 12).Value = 10.25291590449301
$ws.Cells.Item(25, 13).Value = 15.20097497992873
$ws.Cells.Item(25, 14).Value = 21.11724271355905
$ws.Cells.Item(25, 15).Value = 26.76583050629356
